$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.460.50'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.11%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.519.26'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.26%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '615.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.83%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.07%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.518.14'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.15%  '

$ws.Range('E8').Value = '  -0.27%  '

$ws.Range('E9').Value = '  -1.06%  '

$ws.Range('E10').Value = '  -0.79%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.12'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.26%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.426'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.89%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000221'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.40%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.12'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.90%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.114.99'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.29%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.520.24'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.23%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.447.85'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.02%  '

$ws.Range('E18').Value = '  +0.03%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.38'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.24%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.37'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.57%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '444.83'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.60%  '

$ws.Range('E22').Value = '  +1.71%  '

$ws.Range('E23').Value = '  -2.87%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.38'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.62%  '

$ws.Range('E25').Value = '  +6.46%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.660.81'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.44%  '

$ws.Range('E27').Value = '  +0.04%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.21'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.96%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.51'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.25%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.52'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.40%  '

$ws.Range('E31').Value = '  -5.99%  '

$ws.Range('E32').Value = '  +0.10%  '

$ws.Range('E33').Value = '  +4.16%  '

$ws.Range('E34').Value = '  -0.35%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.14'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.85%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.513.21'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.49%  '

$ws.Range('E37').Value = '  -3.30%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.00'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.62%  '

$ws.Range('E39').Value = '  -0.03%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.12%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '177.37'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.56%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0883'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.55%  '

$ws.Range('E43').Value = '  +1.50%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.37%  '

$ws.Range('E45').Value = '  -1.27%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.24'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.57%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '45.00'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.61%  '

$ws.Range('E48').Value = '  -0.22%  '

$ws.Range('E49').Value = '  +1.87%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.60'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.55%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.994'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.74%  '
